# Insert a new weekly price record for "Choclo" (Choclero, Primera)
# reported 2023-01-25 (serial 44951) just above the existing row 127,
# pushing all subsequent records down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 127; everything at/after row 127 shifts
# down by one (this also grows the sheet's used range to row 157).
$ws.Rows.Item(127).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A127").Value = 2
$ws.Range("B127").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C127").Value = "Coquimbo"
$ws.Range("D127").Value = 44951
$ws.Range("E127").Value = 4
$ws.Range("F127").Value = 100112024
$ws.Range("G127").Value = "Choclo"
$ws.Range("H127").Value = "Choclero"
$ws.Range("I127").Value = "Primera"
$ws.Range("J127").Value = 140000
$ws.Range("K127").Value = 200
$ws.Range("L127").Value = 250
$ws.Range("M127").Value = 225
$ws.Range("N127").Value = "$/unidad"
$ws.Range("O127").Value = "Provincia de Limarí"
$ws.Range("P127").Value = 225
$ws.Range("Q127").Value = 1
$ws.Range("R127").Value = "Hortaliza"
